# Insert two new "Memory – LLM Training" slides right after the slide that
# currently ends the "297 / rId13" run (i.e. before the existing "Prompt
# Engineering" slide), pushing it and everything after it two slots later.
# Both new slides reuse the same "Title and Content" layout (slideLayout2)
# that the neighbouring slides already use.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# New slide #1 at position 13: "Memory – LLM Training"
# ---------------------------------------------------------------------
$s1 = $p.Slides.Add(13, 2)

$title1 = $s1.Shapes.Item(1).TextFrame.TextRange
$title1.Text = "Memory – LLM Training"
$title1.ParagraphFormat.Alignment = 2

$body1 = $s1.Shapes.Item(2).TextFrame
$body1.AutoSize = 2
$tr1 = $body1.TextRange
$tr1.Text = "32 bit float – Full precision"
$tr1.InsertAfter("`rModel parameters")
$tr1.InsertAfter("`r4 bytes per parameter")
$tr1.InsertAfter("`rOptimizer")
$tr1.InsertAfter("`r8 bytes per parameter (2 states)")
$tr1.InsertAfter("`rGradients")
$tr1.InsertAfter("`r4 bytes per parameter")
$tr1.InsertAfter("`rActivations ")
$tr1.InsertAfter("`r4 bytes per parameter")
$tr1.InsertAfter("`rTemporary Variables")
$tr1.InsertAfter("`r4 bytes per parameter")
$tr1.InsertAfter("`r")

$body1.TextRange.Paragraphs(3).IndentLevel = 2
$body1.TextRange.Paragraphs(5).IndentLevel = 2
$body1.TextRange.Paragraphs(7).IndentLevel = 2
$body1.TextRange.Paragraphs(9).IndentLevel = 2
$body1.TextRange.Paragraphs(11).IndentLevel = 2
$body1.TextRange.Paragraphs(12).IndentLevel = 2

# ---------------------------------------------------------------------
# New slide #2 at position 14: "Memory – LLM  Training"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(14, 2)

$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
$title2.Text = "Memory – LLM  Training"
$title2.ParagraphFormat.Alignment = 2

$body2 = $s2.Shapes.Item(2).TextFrame
$body2.AutoSize = 2
$tr2 = $body2.TextRange
$tr2.Text = "Store LLM"
$tr2.InsertAfter("`r4 GB memory (1B parameters @ full precision)")
$tr2.InsertAfter("`rTrain LLM")
$tr2.InsertAfter("`r24 GB memory (1B parameters @ full precision)")
$tr2.InsertAfter("`r6 times memory required to store LLM")
$tr2.InsertAfter("`r")

$body2.TextRange.Paragraphs(2).IndentLevel = 2
$body2.TextRange.Paragraphs(4).IndentLevel = 2
$body2.TextRange.Paragraphs(5).IndentLevel = 2
$body2.TextRange.Paragraphs(6).IndentLevel = 2

Write-Output "Slides inserted. Total slide count:"
Write-Output $p.Slides.Count
